$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + date range) ---
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Convert numeric cells to text placeholders ("0" / "***.*") ---
# Use Copy from a stable untouched placeholder cell to preserve style s=14 + shared string.
$ws.Range("G14").Copy($ws.Range("C14"))
$ws.Range("G14").Copy($ws.Range("C22"))
$ws.Range("G14").Copy($ws.Range("D23"))
$ws.Range("H14").Copy($ws.Range("E23"))
$ws.Range("G14").Copy($ws.Range("D26"))
$ws.Range("H14").Copy($ws.Range("E26"))
$ws.Range("G14").Copy($ws.Range("D27"))
$ws.Range("H14").Copy($ws.Range("E27"))
$ws.Range("G14").Copy($ws.Range("C28"))
$ws.Range("G14").Copy($ws.Range("C29"))

# --- Convert text placeholder cells to numeric values ---
$ws.Range("F23").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("F23").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 1

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 83.333333333333
$ws.Range("M15").Value = 37.5
$ws.Range("N15").Value = -52.173913043478
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -20
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = 14.049586776859
$ws.Range("L16").Value = 21.052631578947
$ws.Range("M16").Value = 38
$ws.Range("N16").Value = -75.045207956600
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -21.428571428571
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = 119
$ws.Range("K17").Value = 2.521008403361
$ws.Range("L17").Value = 25.773195876288
$ws.Range("M17").Value = 3.389830508474
$ws.Range("N17").Value = -46.017699115044
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = 33.620689655172
$ws.Range("L18").Value = 3.333333333333
$ws.Range("M18").Value = 72.222222222222
$ws.Range("N18").Value = -72.807017543859
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = -13.333333333333
$ws.Range("I19").Value = 637
$ws.Range("J19").Value = 471
$ws.Range("K19").Value = 35.244161358811
$ws.Range("L19").Value = 32.708333333333
$ws.Range("M19").Value = 10.975609756097
$ws.Range("N19").Value = -14.611260053619
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 9
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 55
$ws.Range("K20").Value = 34.545454545454
$ws.Range("L20").Value = 72.093023255813
$ws.Range("M20").Value = 89.743589743589
$ws.Range("N20").Value = -83.877995642701
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -20.689655172413
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 127
$ws.Range("H21").Value = -9.448818897637
$ws.Range("I21").Value = 1142
$ws.Range("J21").Value = 894
$ws.Range("K21").Value = 27.740492170022
$ws.Range("L21").Value = 27.883538633818
$ws.Range("M21").Value = 22.795698924731
$ws.Range("N21").Value = -55.822050290135
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = 52.631578947368
$ws.Range("L22").Value = 81.25
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 49
$ws.Range("K23").Value = -16.949152542372
$ws.Range("L23").Value = 11.363636363636
$ws.Range("M23").Value = 25.641025641025
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = -10.9375
$ws.Range("I24").Value = 670
$ws.Range("J24").Value = 512
$ws.Range("K24").Value = 30.859375
$ws.Range("L24").Value = 2.290076335877
$ws.Range("M24").Value = -12.646675358539
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -36.842105263157
$ws.Range("I25").Value = 337
$ws.Range("J25").Value = 293
$ws.Range("K25").Value = 15.017064846416
$ws.Range("L25").Value = 51.121076233183
$ws.Range("M25").Value = 15.017064846416
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = 21.428571428571
$ws.Range("L26").Value = 112.5
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 52
$ws.Range("K27").Value = 13.043478260869
$ws.Range("L27").Value = 40.540540540540
$ws.Range("I30").Value = 17
$ws.Range("K30").Value = 466.666666666667
$ws.Range("L30").Value = 21.428571428571
